$wb = $excel.ActiveWorkbook

# "About" sheet: update the "last updated" date in C1 (2024-01-29 -> 2024-04-10)
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = (Get-Date -Year 2024 -Month 4 -Day 10 -Hour 0 -Minute 0 -Second 0).Date

# "MCF" sheet: bump several plant-type capacity factors to 1 (100%)
$wsMCF = $wb.Worksheets.Item("MCF")
$cellsToUpdate = @("B2","B3","B4","B6","B10","B11","B12","B13","B14","B16","B17","B18")
foreach ($addr in $cellsToUpdate) {
    $wsMCF.Range($addr).Value = 1
}

# Update the active selection to match the saved cursor position
$wsMCF.Activate()
$wsMCF.Range("B17").Select()
